$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 6-15 (everything below the new 5-row table)
$ws.Rows("6:15").Delete()

# Update row 2: keep question "вопрос1" in A2, replace answer in B2 with "ответ25"
$ws.Range("B2").Value = "ответ25"

# Row 3: clear the question cell (A3), replace answer with "ответ2"
$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = "ответ2"

# Row 4: keep question "вопрос3" in A4, clear the answer cell (B4)
$ws.Range("B4").ClearContents()

# Row 5: new question/answer pair
$ws.Range("A5").Value = "вопрос 4"
$ws.Range("B5").Value = "ответ 4"

# Update view: zoom to 125% and select B3
$ws.Application.ActiveWindow.Zoom = 125
$ws.Range("B3").Select()
